# Updated cryptos list on Sat Nov 18 06:30:06 UTC 2023 with GitHub Actions
#
# Price (column D) and Volume(1h) (column E) are stored as plain TEXT in
# this sheet, not numbers/percentages - e.g. "241.06" or "  -1.82%  ".
# Writing a number-looking string straight into Range.Value lets Excel's
# automatic type-detection silently convert it to a floating point number,
# which would corrupt values such as "241.06" (two-decimal price) or lose
# the literal formatting. To guarantee the cells stay plain text (matching
# the source data) the new text is first written into a scratch cell
# (prefixed with a leading apostrophe to force text), then copied and
# pasted as values-only into the destination cell - this keeps the
# destination cell's existing style/number format untouched while still
# landing a pure text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

function Set-TextValue($ref, $value) {
    $scratch.Value = "'" + $value
    $scratch.Copy()
    $ws.Range($ref).PasteSpecial(-4163)
}

# Row 2 - Bitcoin
Set-TextValue "D2" "36.360.52"
Set-TextValue "E2" "  -0.25%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.935.39"
Set-TextValue "E3" "  -2.48%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "241.06"
Set-TextValue "E5" "  -1.82%  "

# Row 6 - XRP
Set-TextValue "E6" "  -3.46%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.01%  "

# Row 8 - Solana
Set-TextValue "D8" "56.42"
Set-TextValue "E8" "  -5.76%  "

# Row 9 - Cardano
Set-TextValue "E9" "  -4.95%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0840"
Set-TextValue "E10" "  +1.52%  "

# Row 11 - TRON
Set-TextValue "E11" "  -1.52%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.219.47"
Set-TextValue "E12" "  -2.31%  "

# Row 13 - Polygon
Set-TextValue "E13" "  -8.00%  "

# Row 14 - Avalanche
Set-TextValue "D14" "20.94"
Set-TextValue "E14" "  -11.66%  "

# Row 15 - Chainlink
Set-TextValue "D15" "13.34"
Set-TextValue "E15" "  -5.13%  "

# Row 16 - Polkadot
Set-TextValue "E16" "  -6.68%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.939.26"
Set-TextValue "E17" "  -1.64%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "36.292.17"
Set-TextValue "E18" "  +0.02%  "

# Row 19 - was Litecoin, now ShibaInu
Set-TextValue "B19" "ShibaInu"
Set-TextValue "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.0₃0861"
Set-TextValue "E19" "  -1.77%  "

# Row 20 - was ShibaInu, now Litecoin
Set-TextValue "B20" "Litecoin"
Set-TextValue "C20" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D20" "68.73"
Set-TextValue "E20" "  -2.37%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "226.56"
Set-TextValue "E21" "  -3.49%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.93"
Set-TextValue "E22" "  -7.65%  "

# Row 23 - Dai
Set-TextValue "E23" "  -0.20%  "

# Row 24 - PancakeSwap
Set-TextValue "D24" "2.33"
Set-TextValue "E24" "  -9.95%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.24"
Set-TextValue "E25" "  -2.85%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.22"
Set-TextValue "E26" "  -7.60%  "

# Row 27
Set-TextValue "D27" "160.39"
Set-TextValue "E27" "  -0.94%  "

# Row 28
Set-TextValue "D28" "0.132"
Set-TextValue "E28" "  +0.90%  "

# Row 29
Set-TextValue "D29" "19.09"
Set-TextValue "E29" "  -3.93%  "

# Row 30
Set-TextValue "E30" "  -2.60%  "

# Row 31
Set-TextValue "E31" "  -7.26%  "

# Row 32
Set-TextValue "D32" "4.53"
Set-TextValue "E32" "  -7.97%  "

# Row 33
Set-TextValue "D33" "0.0624"
Set-TextValue "E33" "  -4.19%  "

# Row 34
Set-TextValue "D34" "4.12"
Set-TextValue "E34" "  -7.02%  "

# Row 35
Set-TextValue "E35" "  -0.04%  "

# Row 36
Set-TextValue "D36" "6.03"
Set-TextValue "E36" "  -3.32%  "

# Row 37
Set-TextValue "E37" "  +0.29%  "

# Row 38
Set-TextValue "D38" "2.11"
Set-TextValue "E38" "  -6.76%  "

# Row 39
Set-TextValue "D39" "2.97"
Set-TextValue "E39" "  -2.06%  "

# Row 40
Set-TextValue "D40" "0.0967"
Set-TextValue "E40" "  -0.65%  "

# Row 41
Set-TextValue "E41" "  -1.04%  "

# Row 42
Set-TextValue "D42" "0.0208"
Set-TextValue "E42" "  -3.13%  "

# Row 43
Set-TextValue "E43" "  -8.04%  "

# Row 44
Set-TextValue "D44" "15.40"
Set-TextValue "E44" "  -5.46%  "

# Row 45
Set-TextValue "D45" "1.327.28"
Set-TextValue "E45" "  -3.19%  "

# Row 46
Set-TextValue "E46" "  -7.62%  "

# Row 47
Set-TextValue "D47" "85.40"
Set-TextValue "E47" "  -7.76%  "

# Row 48
Set-TextValue "D48" "7.02"
Set-TextValue "E48" "  -6.45%  "

# Row 49
Set-TextValue "D49" "2.82"
Set-TextValue "E49" "  -0.21%  "

# Row 50
Set-TextValue "D50" "43.65"
Set-TextValue "E50" "  -4.39%  "

# Row 51
Set-TextValue "D51" "2.111.07"
Set-TextValue "E51" "  -2.32%  "

# Clean up the scratch cell used for the text-preserving copy/paste trick.
$scratch.Clear()
